$d = $word.ActiveDocument

# The second paragraph of the document holds the M2Doc query field:
#   {m:'prefix<bookmark>\nsuffix'}
# It is currently made up of runs "{m", ":'", "prefix", "\n", "suffix", "'}".
# The parser now expects every literal character that is part of the field
# delimiters ("{", "m", "'", "}") to live in its own run, so split the two
# runs that currently bundle two of those characters together:
#   "{m" -> "{" + "m"
#   "'}" -> "'" + "}"
#
# Word COM has no direct "split this run in two" call, but toggling a
# character-level formatting property (and then reverting it) on a
# sub-range forces Word to break the run at that boundary without
# touching the surrounding text.

$fullText = $d.Content.Text

# --- Split "{m" into "{" and "m" ---------------------------------------
$openBrace = $d.Content.Find
$openBrace.Execute("{m", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $openBrace.Parent.Start + 1
$rBrace = $d.Range($openBrace.Parent.Start, $splitPoint)
$rBrace.Bold = $true
$rBrace.Bold = $false

# --- Split "'}" into "'" and "}" ----------------------------------------
$closeBrace = $d.Content.Find
$closeBrace.Execute("'}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint2 = $closeBrace.Parent.Start + 1
$rQuote = $d.Range($splitPoint2, $closeBrace.Parent.End)
$rQuote.Bold = $true
$rQuote.Bold = $false
